$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JS-SPA-Self-Evaluation-Protocol")

# Mark "Yes" for the newly implemented town edit/delete (and create) features
$ws.Range("C48").Value = "Yes"
$ws.Range("C49").Value = "Yes"
$ws.Range("C50").Value = "Yes"

# Update the active selection to reflect where the user last worked
$ws.Range("J39").Select()
